# Apply the "additional scraping" update:
#  1. Insert a new "Player Info" sheet at the front of the workbook with
#     player ID/NAME/BATTING_HAND/BOWL_STYLE.
#  2. On the existing "ODI Batting" and "ODI Bowling" sheets, rename the
#     MATCH_CARD_LINK column to MATCH_CODE and replace the full scorecard
#     URL values with just the numeric match code.

$wb = $excel.ActiveWorkbook

function Set-TextValue($cell, $text) {
    # Force the cell to hold a genuine text value (even when the text looks
    # like a number, e.g. "4636") without leaving the cell's visible style
    # changed from its original "Normal" formatting.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# --- 1. New "Player Info" sheet, inserted as the first tab -----------------
# Clone an existing sheet (rather than Worksheets.Add()) so the bold/
# bordered/centred header style used elsewhere in the workbook (style index
# referenced by the other header rows) carries over exactly instead of a
# freshly-synthesised look-alike style.
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Copy($batting, $null)

# NOTE: after Copy(), the *original* `$batting` handle ends up pointing at
# the newly-inserted copy (sheet handles here track by position, and the
# copy is inserted before the source) rather than staying bound to the
# original "ODI Batting" sheet. Re-resolve both sheets fresh, by name/index,
# so the rest of the script edits the correct tabs.
$info = $wb.Worksheets.Item(1)
$info.Name = "Player Info"
$batting = $wb.Worksheets.Item("ODI Batting")

$info.Range("A1").Value = "ID"
$info.Range("B1").Value = "NAME"
$info.Range("C1").Value = "BATTING_HAND"
$info.Range("D1").Value = "BOWL_STYLE"
$info.Range("E1:J1").Clear()

Set-TextValue $info.Range("A2") "6854"
$info.Range("B2").Value = "Yannic Cariah"
$info.Range("C2").Value = "Left Handed"
$info.Range("D2").Value = "Right Arm Leg Break"
$info.Range("E2:J6").Clear()
$info.Range("A3:D6").Clear()

# --- 2. ODI Batting: MATCH_CARD_LINK -> MATCH_CODE --------------------------
$batting.Range("D1").Value = "MATCH_CODE"

$battingCodes = @("4636", "4639", "4642", "4727", "4731")
for ($i = 0; $i -lt $battingCodes.Length; $i++) {
    $row = $i + 2
    Set-TextValue $batting.Cells.Item($row, 4) $battingCodes[$i]
}

# --- 3. ODI Bowling: MATCH_CARD_LINK -> MATCH_CODE --------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowling.Range("B1").Value = "MATCH_CODE"

$bowlingCodes = @("4636", "4639", "4642", "4727", "4731")
for ($i = 0; $i -lt $bowlingCodes.Length; $i++) {
    $row = $i + 2
    Set-TextValue $bowling.Cells.Item($row, 2) $bowlingCodes[$i]
}
